$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: update existing K6 value
$ws.Range("K6").Value = 0.5403464745801891

# Row 7: update existing J7 value, add new K7 value
$ws.Range("J7").Value = 0.4969491838668565
$ws.Range("K7").Value = 0.2970525035592049

# Row 8: update existing I8 value, add new J8 value
$ws.Range("I8").Value = 0.5011245128056051
$ws.Range("J8").Value = 0.2858677898194339

# Row 9: update existing H9 value, add new I9 value
$ws.Range("H9").Value = 0.4852787037784192
$ws.Range("I9").Value = 0.2775335613519331

# Row 10: update existing G10 value, add new H10 value
$ws.Range("G10").Value = 0.4539510573947921
$ws.Range("H10").Value = 0.2743085116504074

# Row 11: update existing F11 value, add new G11 value
$ws.Range("F11").Value = 0.4663391832225094
$ws.Range("G11").Value = 0.2534447081011285

# Row 12: update existing E12 value, add new F12 value
$ws.Range("E12").Value = 0.4814444548743619
$ws.Range("F12").Value = 0.2766837437271186

# Row 13: update existing D13 value, add new E13 value
$ws.Range("D13").Value = 0.4184715358843989
$ws.Range("E13").Value = 0.2867219094086165

# Row 14: update existing C14 value, add new D14 value
$ws.Range("C14").Value = 0.5177895860664353
$ws.Range("D14").Value = 0.1751453671933744

# Row 15: update existing B15 value, add new C15 value
$ws.Range("B15").Value = 0.5618492773058843
$ws.Range("C15").Value = 0.1965658720679752

# Row 16: add new B16 value
$ws.Range("B16").Value = 0.4328090033804217
